# Fix exploration time formula (#13)
# Updates Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M), and Std_Experiment_Time (N) columns
# for rows 2-13 to reflect the corrected formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 6.554423839999999
$ws.Range("H2").Value2 = 357.76253134
$ws.Range("M2").Value2 = 0.496689685637715
$ws.Range("N2").Value2 = 46.44027956847143
$ws.Range("G3").Value2 = 6.93212394
$ws.Range("H3").Value2 = 633.4908791399999
$ws.Range("M3").Value2 = 0.689183557178817
$ws.Range("N3").Value2 = 128.3414278038793
$ws.Range("G4").Value2 = 2.69610748
$ws.Range("H4").Value2 = 85.93855697999999
$ws.Range("M4").Value2 = 0.4394778409328411
$ws.Range("N4").Value2 = 24.90319498582187
$ws.Range("G5").Value2 = 2.806555679999999
$ws.Range("H5").Value2 = 138.01340042
$ws.Range("M5").Value2 = 0.3700454724131137
$ws.Range("N5").Value2 = 34.83768979095407
$ws.Range("G6").Value2 = 0.9318887
$ws.Range("H6").Value2 = 15.4618549
$ws.Range("M6").Value2 = 0.2470606015498756
$ws.Range("N6").Value2 = 7.274285970939513
$ws.Range("G7").Value2 = 1.09178768
$ws.Range("H7").Value2 = 28.71313114
$ws.Range("M7").Value2 = 0.2378947320258877
$ws.Range("N7").Value2 = 11.88142715302962
$ws.Range("G8").Value2 = 0.46838378
$ws.Range("H8").Value2 = 5.01634664
$ws.Range("M8").Value2 = 0.1619809555640492
$ws.Range("N8").Value2 = 3.004317489136955
$ws.Range("G9").Value2 = 0.54759462
$ws.Range("H9").Value2 = 9.974282200000001
$ws.Range("M9").Value2 = 0.1531322757743228
$ws.Range("N9").Value2 = 5.775106211956548
$ws.Range("G10").Value2 = 0.2488806
$ws.Range("H10").Value2 = 1.86258992
$ws.Range("M10").Value2 = 0.1057977267790913
$ws.Range("N10").Value2 = 1.183603468563346
$ws.Range("G11").Value2 = 0.31358608
$ws.Range("H11").Value2 = 4.444253300000001
$ws.Range("M11").Value2 = 0.1076541527514981
$ws.Range("N11").Value2 = 3.013857900209201
$ws.Range("G12").Value2 = 0.16727788
$ws.Range("H12").Value2 = 1.07234594
$ws.Range("M12").Value2 = 0.08858614315921678
$ws.Range("N12").Value2 = 0.8449439894544175
$ws.Range("G13").Value2 = 0.1795166
$ws.Range("H13").Value2 = 2.10835274
$ws.Range("M13").Value2 = 0.06813177213303626
$ws.Range("N13").Value2 = 1.58756538453717
